$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Poultry" row entirely (row 11); all rows below shift up by one.
$ws.Rows.Item(11).Delete()

# Rename the header for column A.
$ws.Range("A1").Value = "Main Category"

# Recomputed nutrition statistics for each remaining category (rows 2-14).
# Row 2: Beverages
$ws.Range("B2").Value = 0.7868965517241379
$ws.Range("C2").Value = 0.5206896551724138
$ws.Range("D2").Value = 9.204137931034483
$ws.Range("E2").Value = 87.32638888888889
$ws.Range("F2").Value = 26.5
$ws.Range("G2").Value = 8.106206896551726
$ws.Range("H2").Value = 58.7448275862069
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 12.67361111111111

# Row 3: Condiments and sauces
$ws.Range("B3").Value = 11.17368421052632
$ws.Range("C3").Value = 8.25
$ws.Range("D3").Value = 25.12105263157895
$ws.Range("E3").Value = 33.05263157894737
$ws.Range("F3").Value = 249.3684210526316
$ws.Range("G3").Value = 6.850000000000001
$ws.Range("H3").Value = 264.921052631579
$ws.Range("I3").Value = 94.31578947368421
$ws.Range("J3").Value = 66.94736842105263

# Row 4: Dairy and alternatives
$ws.Range("B4").Value = 11.74583333333333
$ws.Range("C4").Value = 14.54027777777778
$ws.Range("D4").Value = 5.888888888888889
$ws.Range("E4").Value = 67.6875
$ws.Range("F4").Value = 284.8333333333333
$ws.Range("G4").Value = 5.770138888888889
$ws.Range("H4").Value = 201.6597222222222
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 32.3125

# Row 5: Eggs
$ws.Range("B5").Value = 12.56
$ws.Range("C5").Value = 13.23333333333333
$ws.Range("D5").Value = 1.926666666666667
$ws.Range("E5").Value = 72.33333333333333
$ws.Range("F5").Value = 211.8666666666667
$ws.Range("G5").Value = 1.92
$ws.Range("H5").Value = 177.0666666666667
$ws.Range("I5").Value = 96.8
$ws.Range("J5").Value = 27.66666666666667

# Row 6: Fish
$ws.Range("B6").Value = 15.62995594713656
$ws.Range("C6").Value = 7.644493392070484
$ws.Range("D6").Value = 4.531277533039647
$ws.Range("E6").Value = 71.66964285714286
$ws.Range("F6").Value = 194.4096916299559
$ws.Range("G6").Value = 1.165198237885463
$ws.Range("H6").Value = 150.1013215859031
$ws.Range("I6").Value = 88.98360655737704
$ws.Range("J6").Value = 28.33035714285714

# Row 7: Fruits and nuts
$ws.Range("B7").Value = 4.614364640883978
$ws.Range("C7").Value = 10.12099447513812
$ws.Range("D7").Value = 19.88397790055249
$ws.Range("E7").Value = 60.26519337016575
$ws.Range("F7").Value = 121.6850828729282
$ws.Range("G7").Value = 16.2939226519337
$ws.Range("H7").Value = 199.1436464088398
$ws.Range("I7").Value = 90.44692737430168
$ws.Range("J7").Value = 39.73480662983425

# Row 8: Grains and cereals
$ws.Range("B8").Value = 7.895876288659793
$ws.Range("C8").Value = 9.149484536082474
$ws.Range("D8").Value = 49.45541237113402
$ws.Range("E8").Value = 28.59536082474227
$ws.Range("F8").Value = 182.9716494845361
$ws.Range("G8").Value = 11.57268041237113
$ws.Range("H8").Value = 321.6417525773196
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 71.40463917525773

# Row 9: Legumes
$ws.Range("B9").Value = 13.90645161290323
$ws.Range("C9").Value = 1.990322580645161
$ws.Range("D9").Value = 24.88387096774193
$ws.Range("E9").Value = 50.54838709677419
$ws.Range("F9").Value = 222.3225806451613
$ws.Range("G9").Value = 2.7
$ws.Range("H9").Value = 190.6451612903226
$ws.Range("I9").Value = 99.64516129032258
$ws.Range("J9").Value = 49.45161290322581

# Row 10: Miscellaneous
$ws.Range("B10").Value = 6.017021276595744
$ws.Range("C10").Value = 22.48662613981763
$ws.Range("D10").Value = 13.8258358662614
$ws.Range("E10").Value = 55.55487804878049
$ws.Range("F10").Value = 127.0501567398119
$ws.Range("G10").Value = 4.196625766871166
$ws.Range("H10").Value = 288.7055214723927
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 44.44512195121951

# Row 11: Red meat
$ws.Range("B11").Value = 18.39186046511628
$ws.Range("C11").Value = 11.73875968992248
$ws.Range("D11").Value = 2.857364341085272
$ws.Range("E11").Value = 66.67829457364341
$ws.Range("F11").Value = 196.9182879377432
$ws.Range("G11").Value = 0.6833333333333333
$ws.Range("H11").Value = 191.3759689922481
$ws.Range("I11").Value = 97.33070866141732
$ws.Range("J11").Value = 33.32170542635659

# Row 12: Starchy vegetables
$ws.Range("B12").Value = 2.323529411764706
$ws.Range("C12").Value = 5.8
$ws.Range("D12").Value = 18.73529411764706
$ws.Range("E12").Value = 70.47058823529412
$ws.Range("F12").Value = 72.41176470588235
$ws.Range("G12").Value = 2.170588235294117
$ws.Range("H12").Value = 140.4117647058823
$ws.Range("I12").Value = 96.11764705882354
$ws.Range("J12").Value = 29.52941176470588

# Row 13: Sweets and snacks
$ws.Range("B13").Value = 3.512903225806452
$ws.Range("C13").Value = 9.835483870967741
$ws.Range("D13").Value = 20.18064516129032
$ws.Range("E13").Value = 64.83870967741936
$ws.Range("F13").Value = 102.6451612903226
$ws.Range("G13").Value = 18.70645161290323
$ws.Range("H13").Value = 186.7741935483871
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 35.16129032258065

# Row 14: Vegetables
$ws.Range("B14").Value = 2.003164556962025
$ws.Range("C14").Value = 1.654430379746836
$ws.Range("D14").Value = 5.648101265822785
$ws.Range("E14").Value = 87.91082802547771
$ws.Range("F14").Value = 50.28481012658228
$ws.Range("G14").Value = 3.557961783439491
$ws.Range("H14").Value = 51.03797468354431
$ws.Range("I14").Value = 89.52903225806452
$ws.Range("J14").Value = 12.08917197452229
